# Adapt column header formatting to the respective input file names:
#   "<name>_old" -> "<name>_FV2304"
#   "<name>_new" -> "<name>_FV2310"
# (column "diff" is untouched), then wrap the sheet's data range in an
# Excel Table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastCol = $usedRange.Columns.Count
$lastRow = $usedRange.Rows.Count

# Rename header row (row 1) suffixes based on each header's current text.
for ($col = 1; $col -le $lastCol; $col++) {
    $headerCell = $ws.Cells.Item(1, $col)
    $headerText = $headerCell.Value2
    if ($headerText -like "*_old") {
        $headerCell.Value2 = $headerText -replace "_old$", "_FV2304"
    } elseif ($headerText -like "*_new") {
        $headerCell.Value2 = $headerText -replace "_new$", "_FV2310"
    }
}

# Wrap the data range in an Excel Table ("Table1") with a header row.
$dataRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$lo = $ws.ListObjects.Add(1, $dataRange, [System.Type]::Missing, 1)
$lo.Name = "Table1"

# Freeze the header row (View > Freeze Panes > Freeze Top Row).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
